# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (Total) sheet,
#    populated with the quarterly fund-holding detail rows.
# 2. Prepend a corresponding summary row to the "总计" sheet and renumber
#    the existing index column.

$wb = $excel.ActiveWorkbook

$sourceDetail = $wb.Worksheets.Item("2021-Q4")
$totalSheet   = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------------
# 1. New "2022-Q1" detail sheet (inserted immediately before "总计")
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# NOTE: adding a sheet silently re-targets any worksheet object obtained
# *before* the Add() call onto whatever sheet is now active, so every
# worksheet reference we still need has to be re-fetched afterwards.
$sourceDetail = $wb.Worksheets.Item("2021-Q4")
$totalSheet   = $wb.Worksheets.Item("总计")

# Match the page margins used by the sibling quarterly sheets (inches, as
# COM margins are expressed in points: 0.75in=54pt, 1in=72pt, 0.5in=36pt).
$newSheet.PageSetup.LeftMargin   = 54
$newSheet.PageSetup.RightMargin  = 54
$newSheet.PageSetup.TopMargin    = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Pull over the header-row and index-column styling from an existing
# quarterly sheet so the new sheet matches the established look.
$sourceDetail.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$sourceDetail.Range("A2:A12").Copy()
$newSheet.Range("A2:A12").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Fund code (B) and the numeric-looking metrics (D:G) must stay text so
# leading zeros / original formatting are preserved instead of Excel
# auto-coercing them into numbers.
$newSheet.Range("B2:B12").NumberFormat = "@"
$newSheet.Range("D2:G12").NumberFormat = "@"

function Set-FundRow {
    param($Row, $Idx, $Code, $FundName, $Scale, $StockPos, $PosPct, $MarketValue, $Rank)
    $newSheet.Cells.Item($Row, 1).Value = $Idx
    $newSheet.Cells.Item($Row, 2).Value = $Code
    $newSheet.Cells.Item($Row, 3).Value = $FundName
    $newSheet.Cells.Item($Row, 4).Value = $Scale
    $newSheet.Cells.Item($Row, 5).Value = $StockPos
    $newSheet.Cells.Item($Row, 6).Value = $PosPct
    $newSheet.Cells.Item($Row, 7).Value = $MarketValue
    $newSheet.Cells.Item($Row, 8).Value = $Rank
}

Set-FundRow 2  0  "516150" "嘉实中证稀土产业ETF"                         "25.17" "99.75" "3.41" "0.8583" 10
Set-FundRow 3  1  "516780" "华泰柏瑞中证稀土产业ETF"                     "11.06" "98.70" "3.41" "0.3771" 10
Set-FundRow 4  2  "006973" "太平睿盈混合A"                               "8.69"  "29.33" "1.94" "0.1686" 2
Set-FundRow 5  3  "011346" "淳厚鑫淳一年持有期混合型证券投资基金"         "5.75"  "67.80" "2.65" "0.1524" 8
Set-FundRow 6  4  "159715" "易方达中证稀土产业ETF"                       "3.42"  "99.06" "3.38" "0.1156" 10
Set-FundRow 7  5  "012454" "淳厚鑫悦混合A"                               "3.31"  "76.84" "3.49" "0.1155" 3
Set-FundRow 8  6  "159713" "富国中证稀土产业交易型开放式指数证券投资基金" "3.26"  "99.26" "3.39" "0.1105" 10
Set-FundRow 9  7  "007669" "太平睿盈混合C"                               "2.21"  "29.33" "1.94" "0.0429" 2
Set-FundRow 10 8  "012455" "淳厚鑫悦混合C"                               "0.79"  "76.84" "3.49" "0.0276" 3
Set-FundRow 11 9  "160639" "鹏华中证高铁产业指数（LOF）"                 "0.89"  "94.72" "2.75" "0.0245" 6
Set-FundRow 12 10 "002161" "银华万物互联灵活配置混合"                    "1.06"  "20.06" "1.34" "0.0142" 6

# ---------------------------------------------------------------------------
# 2. Update the "总计" (Total) sheet: insert the 2022-Q1 summary row on top
#    and renumber the existing rows' index column.
# ---------------------------------------------------------------------------
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A7").PasteSpecial(-4122)

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 11
$totalSheet.Cells.Item(2, 4).Value = 2.01

$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(3, 2).Value = "2021-Q4"
$totalSheet.Cells.Item(3, 3).Value = 12
$totalSheet.Cells.Item(3, 4).Value = 2.1

$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(4, 2).Value = "2021-Q3"
$totalSheet.Cells.Item(4, 3).Value = 4
$totalSheet.Cells.Item(4, 4).Value = 0.18

$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(5, 2).Value = "2021-Q2"
$totalSheet.Cells.Item(5, 3).Value = 4
$totalSheet.Cells.Item(5, 4).Value = 1.02

$totalSheet.Cells.Item(6, 1).Value = 4
$totalSheet.Cells.Item(6, 2).Value = "2021-Q1"
$totalSheet.Cells.Item(6, 3).Value = 2
$totalSheet.Cells.Item(6, 4).Value = 0.1

$totalSheet.Cells.Item(7, 1).Value = 5
$totalSheet.Cells.Item(7, 2).Value = "2020-Q4"
$totalSheet.Cells.Item(7, 3).Value = 2
$totalSheet.Cells.Item(7, 4).Value = 0.12

# Restore the originally active sheet/selection.
$wb.Worksheets.Item("2020-Q4").Select()
